# svmvHomiEdad.xlsx: add the 2023 column of data (SNMV homicide-by-age
# table), pushing the existing "Total" column one slot to the right and
# refreshing the row/column totals so they include the new year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H. This shifts the old "Total" column (H) to I,
# carrying its formulas/values/styles along, and leaves a blank H ready
# for the new 2023 data.
$ws.Columns("H").Insert()

# New column header.
$ws.Range("H1").Value2 = 2023

# New 2023 counts per age group (rows 2-13).
$h_values = @(3, 1, 2, 6, 3, 6, 5, 5, 3, 1, 11, 0)
for ($i = 0; $i -lt $h_values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value2 = $h_values[$i]
}

# Row totals (column I) now need to sum through the new H column instead
# of stopping at G.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 9).Formula = "=SUM(B${row}:H${row})"
}

# Column totals (row 14): H14 becomes a plain entered total for the new
# 2023 column, and I14 becomes the grand-total column's sum.
$ws.Range("H14").Value2 = 46
$ws.Range("I14").Formula = "=SUM(I2:I13)"

# Leave the selection where the editor apparently left off.
$ws.Range("I11").Select() | Out-Null
